# Auto-generated: update Bundesliga dataset stats (goals/xG columns) for rows 2-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 12
$ws.Cells.Item(2, 3).Value = 1.71
$ws.Cells.Item(2, 4).Value = 7.4
$ws.Cells.Item(2, 5).Value = 6.7
$ws.Cells.Item(2, 6).Value = 12.1
$ws.Cells.Item(2, 7).Value = 1.06
$ws.Cells.Item(2, 8).Value = 1.84
$ws.Cells.Item(2, 9).Value = 0.95
$ws.Cells.Item(2, 10).Value = 1.73

# Row 3
$ws.Cells.Item(3, 2).Value = 27
$ws.Cells.Item(3, 3).Value = 3.86
$ws.Cells.Item(3, 4).Value = 18.8
$ws.Cells.Item(3, 5).Value = 15.8
$ws.Cells.Item(3, 6).Value = 28.3
$ws.Cells.Item(3, 7).Value = 2.69
$ws.Cells.Item(3, 8).Value = 4.47
$ws.Cells.Item(3, 9).Value = 2.26
$ws.Cells.Item(3, 10).Value = 4.05

# Row 4
$ws.Cells.Item(4, 2).Value = 13
$ws.Cells.Item(4, 3).Value = 1.86
$ws.Cells.Item(4, 4).Value = 11.7
$ws.Cells.Item(4, 5).Value = 10.9
$ws.Cells.Item(4, 6).Value = 20.4
$ws.Cells.Item(4, 7).Value = 1.67
$ws.Cells.Item(4, 8).Value = 3.03
$ws.Cells.Item(4, 9).Value = 1.56

# Row 5
$ws.Cells.Item(5, 2).Value = 19
$ws.Cells.Item(5, 3).Value = 2.71
$ws.Cells.Item(5, 4).Value = 10.9
$ws.Cells.Item(5, 5).Value = 10.1
$ws.Cells.Item(5, 6).Value = 18.9
$ws.Cells.Item(5, 7).Value = 1.55
$ws.Cells.Item(5, 8).Value = 2.81
$ws.Cells.Item(5, 9).Value = 1.44
$ws.Cells.Item(5, 10).Value = 2.69

# Row 6
$ws.Cells.Item(6, 2).Value = 10
$ws.Cells.Item(6, 3).Value = 1.43
$ws.Cells.Item(6, 4).Value = 10.3
$ws.Cells.Item(6, 5).Value = 7.9
$ws.Cells.Item(6, 6).Value = 14.6
$ws.Cells.Item(6, 7).Value = 1.47
$ws.Cells.Item(6, 8).Value = 2.42

# Row 7
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 0.86
$ws.Cells.Item(7, 4).Value = 8.1
$ws.Cells.Item(7, 5).Value = 8.1
$ws.Cells.Item(7, 6).Value = 14.5
$ws.Cells.Item(7, 7).Value = 1.16
$ws.Cells.Item(7, 8).Value = 2.07
$ws.Cells.Item(7, 9).Value = 1.16
$ws.Cells.Item(7, 10).Value = 2.07

# Row 8
$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(8, 4).Value = 7.9
$ws.Cells.Item(8, 5).Value = 7.9
$ws.Cells.Item(8, 6).Value = 14.3
$ws.Cells.Item(8, 7).Value = 1.13
$ws.Cells.Item(8, 8).Value = 2.05
$ws.Cells.Item(8, 9).Value = 1.13
$ws.Cells.Item(8, 10).Value = 2.05

# Row 9
$ws.Cells.Item(9, 2).Value = 6
$ws.Cells.Item(9, 3).Value = 0.86
$ws.Cells.Item(9, 4).Value = 9.6
$ws.Cells.Item(9, 5).Value = 9.6
$ws.Cells.Item(9, 6).Value = 17.4
$ws.Cells.Item(9, 7).Value = 1.37
$ws.Cells.Item(9, 8).Value = 2.49
$ws.Cells.Item(9, 9).Value = 1.37
$ws.Cells.Item(9, 10).Value = 2.49

# Row 10
$ws.Cells.Item(10, 2).Value = 12
$ws.Cells.Item(10, 3).Value = 1.71
$ws.Cells.Item(10, 4).Value = 11.3
$ws.Cells.Item(10, 5).Value = 9.699999999999999
$ws.Cells.Item(10, 6).Value = 18.3
$ws.Cells.Item(10, 7).Value = 1.62
$ws.Cells.Item(10, 8).Value = 2.83
$ws.Cells.Item(10, 9).Value = 1.39
$ws.Cells.Item(10, 10).Value = 2.61

# Row 11
$ws.Cells.Item(11, 2).Value = 12
$ws.Cells.Item(11, 3).Value = 1.71
$ws.Cells.Item(11, 4).Value = 9.6
$ws.Cells.Item(11, 5).Value = 9.6
$ws.Cells.Item(11, 6).Value = 17.5
$ws.Cells.Item(11, 7).Value = 1.37
$ws.Cells.Item(11, 8).Value = 2.5
$ws.Cells.Item(11, 9).Value = 1.37
$ws.Cells.Item(11, 10).Value = 2.5

# Row 12
$ws.Cells.Item(12, 2).Value = 16
$ws.Cells.Item(12, 3).Value = 2.29
$ws.Cells.Item(12, 4).Value = 13.4
$ws.Cells.Item(12, 5).Value = 11.1
$ws.Cells.Item(12, 6).Value = 18.7
$ws.Cells.Item(12, 7).Value = 1.92
$ws.Cells.Item(12, 8).Value = 3
$ws.Cells.Item(12, 9).Value = 1.59
$ws.Cells.Item(12, 10).Value = 2.67

# Row 13
$ws.Cells.Item(13, 2).Value = 8
$ws.Cells.Item(13, 3).Value = 1.14
$ws.Cells.Item(13, 4).Value = 9.699999999999999
$ws.Cells.Item(13, 5).Value = 8.199999999999999
$ws.Cells.Item(13, 6).Value = 14.3
$ws.Cells.Item(13, 7).Value = 1.39
$ws.Cells.Item(13, 8).Value = 2.27
$ws.Cells.Item(13, 9).Value = 1.17
$ws.Cells.Item(13, 10).Value = 2.05

# Row 14
$ws.Cells.Item(14, 2).Value = 10
$ws.Cells.Item(14, 3).Value = 1.43
$ws.Cells.Item(14, 4).Value = 12.4
$ws.Cells.Item(14, 5).Value = 11.6
$ws.Cells.Item(14, 6).Value = 19.1
$ws.Cells.Item(14, 7).Value = 1.77
$ws.Cells.Item(14, 8).Value = 2.85
$ws.Cells.Item(14, 9).Value = 1.66
$ws.Cells.Item(14, 10).Value = 2.74

# Row 15
$ws.Cells.Item(15, 3).Value = 1.14
$ws.Cells.Item(15, 4).Value = 7.9
$ws.Cells.Item(15, 5).Value = 6.3
$ws.Cells.Item(15, 6).Value = 11.2
$ws.Cells.Item(15, 7).Value = 1.13
$ws.Cells.Item(15, 8).Value = 1.83
$ws.Cells.Item(15, 10).Value = 1.6

# Row 16
$ws.Cells.Item(16, 2).Value = 11
$ws.Cells.Item(16, 3).Value = 1.57
$ws.Cells.Item(16, 4).Value = 9.699999999999999
$ws.Cells.Item(16, 5).Value = 8.1
$ws.Cells.Item(16, 6).Value = 15.5
$ws.Cells.Item(16, 7).Value = 1.39
$ws.Cells.Item(16, 8).Value = 2.44
$ws.Cells.Item(16, 9).Value = 1.16
$ws.Cells.Item(16, 10).Value = 2.22

# Row 17
$ws.Cells.Item(17, 2).Value = 11
$ws.Cells.Item(17, 3).Value = 1.57
$ws.Cells.Item(17, 4).Value = 9.1
$ws.Cells.Item(17, 5).Value = 8.300000000000001
$ws.Cells.Item(17, 6).Value = 13.7
$ws.Cells.Item(17, 7).Value = 1.3
$ws.Cells.Item(17, 8).Value = 2.07
$ws.Cells.Item(17, 9).Value = 1.19
$ws.Cells.Item(17, 10).Value = 1.96

# Row 18
$ws.Cells.Item(18, 2).Value = 11
$ws.Cells.Item(18, 3).Value = 1.57
$ws.Cells.Item(18, 4).Value = 9.699999999999999
$ws.Cells.Item(18, 5).Value = 7.4
$ws.Cells.Item(18, 6).Value = 13.5
$ws.Cells.Item(18, 7).Value = 1.39
$ws.Cells.Item(18, 8).Value = 2.26
$ws.Cells.Item(18, 10).Value = 1.92

# Row 19
$ws.Cells.Item(19, 3).Value = 1.14
$ws.Cells.Item(19, 4).Value = 9.9
$ws.Cells.Item(19, 5).Value = 9.199999999999999
$ws.Cells.Item(19, 6).Value = 16.9
$ws.Cells.Item(19, 7).Value = 1.42
$ws.Cells.Item(19, 8).Value = 2.53
$ws.Cells.Item(19, 9).Value = 1.31
$ws.Cells.Item(19, 10).Value = 2.42

Write-Host "Updated Bundesliga stats for rows 2-19"
